$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the "Source" / route column header
$ws.Range("G1").Value = "Source"

# Translate "SMD capacitor" rows (component type)
$ws.Range("B2").Value = "SMD capacitor"
$ws.Range("B3").Value = "SMD capacitor"
$ws.Range("B4").Value = "SMD capacitor"

# Translate "SMD resistor" rows (component type)
$ws.Range("B7").Value = "SMD resistor"
$ws.Range("B8").Value = "SMD resistor"
$ws.Range("B9").Value = "SMD resistor"

# Translate supplier/source values
$ws.Range("G2").Value = "LCSC"
$ws.Range("G3").Value = "LCSC"
$ws.Range("G4").Value = "LCSC"
$ws.Range("G5").Value = "LCSC"
$ws.Range("G6").Value = "LCSC"
$ws.Range("G7").Value = "LCSC"
$ws.Range("G8").Value = "LCSC"
$ws.Range("G9").Value = "LCSC"
$ws.Range("G10").Value = "Taobao"
